$wb = $excel.ActiveWorkbook

# Mapping of sheet name -> row -> new value for column F ("想去人数" / want-to-go count)
$updates = @{
    "展览" = @{
        4 = 6439
        5 = 2565
        6 = 662
        7 = 91
        8 = 3159
        10 = 162
        12 = 7836
        13 = 398
        14 = 72
        15 = 111
        16 = 21
        18 = 269
        19 = 122
        20 = 35
        21 = 274
        22 = 9652
        23 = 27
        24 = 273
        25 = 22
        26 = 29
        27 = 361
        28 = 137
        32 = 76
        33 = 2631
        36 = 19
        37 = 1496
        38 = 815
        39 = 3985
        40 = 228
        41 = 926
        42 = 1204
        43 = 128
        44 = 270
        45 = 120
        46 = 22
        48 = 48
        49 = 71
        50 = 25
    }
    "演出" = @{
        3 = 2
        9 = 31
        23 = 13
    }
    "全部类型" = @{
        6 = 6439
        7 = 2565
        8 = 662
        9 = 91
        10 = 3159
        12 = 31
        14 = 162
        16 = 7836
        17 = 398
        18 = 72
        19 = 111
        20 = 269
        21 = 35
        22 = 274
        23 = 9652
        24 = 273
        25 = 29
        26 = 361
        27 = 137
        30 = 76
        31 = 2631
        33 = 19
        34 = 1496
        35 = 815
        37 = 3985
        38 = 229
        39 = 926
        41 = 1204
        42 = 128
        43 = 270
        44 = 13
        45 = 120
        46 = 22
        48 = 48
        49 = 71
        50 = 25
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $updates[$sheetName]
    foreach ($row in $rowMap.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowMap[$row]
    }
}

Write-Output "Done updating F column values across 展览/演出/全部类型 sheets"